$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Week 5" is already the active sheet (tabSelected=1 / activeTab=4 in the
# source file), but resolve it explicitly too, just to be safe.
$ws = $wb.Worksheets.Item("Week 5")
$ws.Activate() | Out-Null

# New timesheet entry for 2/8/2018 (serial 41677 under the workbook's
# 1904 date system): pull the date/time number formats from row 2 via a
# copy + paste-special (formats only) so the shared style entries for the
# date (s=5) and time (s=6) columns are reused rather than duplicated.
$ws.Range("A2:C2").Copy() | Out-Null
$ws.Range("A3:C3").PasteSpecial(-4122) | Out-Null

$ws.Range("A3").Value = 41677
$ws.Range("B3").Value = 0.41666666666666669
$ws.Range("C3").Value = 0.5

# Description cell: wrap text like the other multi-line description cells
# (reuses the existing wrap-text style rather than creating a new one).
$ws.Range("D3").WrapText = $true
$ws.Range("D3").Value = "Created sql database dump for alchemortem, created master layout for laravel application"

$ws.Range("E3").Value = 2

# Move the selection the way the author's session ended up.
$ws.Range("E4").Select() | Out-Null
